$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price strings such as "573.68" or "64.281.92".
# Excel auto-detects numeric-looking text and would silently convert
# single-dot values (e.g. "573.50") into numbers, losing the original
# text formatting. Force the column to Text first, write the values,
# then restore the default "Normal" style so no stray number format is
# left behind on the cells (matches the source workbook, which uses the
# default style for these cells).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "64.228.82"
$ws.Range("E2").Value = "  -1.23%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.439.10"
$ws.Range("E3").Value = "  -0.33%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "573.50"
$ws.Range("E5").Value = "  -0.43%  "

# Row 6 - Solana
$ws.Range("D6").Value = "164.60"
$ws.Range("E6").Value = "  +2.95%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.05%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.438.37"
$ws.Range("E8").Value = "  -0.37%  "

# Row 9 - XRP
$ws.Range("D9").Value = "0.554"
$ws.Range("E9").Value = "  -6.11%  "

# Row 10 - Toncoin
$ws.Range("D10").Value = "7.29"
$ws.Range("E10").Value = "  +0.58%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -2.26%  "

# Row 12 - Cardano
$ws.Range("D12").Value = "0.426"
$ws.Range("E12").Value = "  -5.08%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "4.032.06"
$ws.Range("E13").Value = "  -0.29%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +1.33%  "

# Row 15 - Avalanche
$ws.Range("D15").Value = "27.30"
$ws.Range("E15").Value = "  -1.89%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -6.98%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "64.280.66"
$ws.Range("E17").Value = "  -1.17%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.504.79"
$ws.Range("E18").Value = "  +1.43%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  -3.52%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "13.71"
$ws.Range("E20").Value = "  -1.78%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "379.18"
$ws.Range("E21").Value = "  -1.11%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "7.82"
$ws.Range("E22").Value = "  -2.37%  "

# Row 23 - Dai
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.30%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "71.42"
$ws.Range("E24").Value = "  -1.16%  "

# Row 25 - Polygon
$ws.Range("D25").Value = "0.521"
$ws.Range("E25").Value = "  -5.56%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  -1.89%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").Value = "9.59"
$ws.Range("E27").Value = "  -3.20%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +0.03%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  -0.01%  "

# Row 30 - NEARProtocol
$ws.Range("D30").Value = "6.10"
$ws.Range("E30").Value = "  -0.87%  "

# Row 31 - Fetch.AI
$ws.Range("E31").Value = "  -5.73%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  -0.24%  "

# Row 33 - EthereumClassic
$ws.Range("D33").Value = "23.02"
$ws.Range("E33").Value = "  -1.45%  "

# Row 34 - Aptos
$ws.Range("E34").Value = "  +0.81%  "

# Row 35 - ImmutableX
$ws.Range("E35").Value = "  -4.38%  "

# Row 36 - Monero
$ws.Range("D36").Value = "159.99"
$ws.Range("E36").Value = "  -0.63%  "

# Row 37 - Mantle
$ws.Range("D37").Value = "0.863"
$ws.Range("E37").Value = "  +10.90%  "

# Row 38 - Stacks
$ws.Range("E38").Value = "  -4.67%  "

# Row 39 - Maker
$ws.Range("D39").Value = "2.816.63"
$ws.Range("E39").Value = "  -3.28%  "

# Row 40 - now Hedera (was EnergySwap)
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.0730"
$ws.Range("E40").Value = "  -3.24%  "

# Row 41 - now EnergySwap (was Hedera)
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "26.08"
$ws.Range("E41").Value = "  -1.63%  "

# Row 42 - now OKB (was InjectiveProtocol)
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "42.97"
$ws.Range("E42").Value = "  -0.59%  "

# Row 43 - now InjectiveProtocol (was OKB)
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "26.46"
$ws.Range("E43").Value = "  +0.77%  "

# Row 44 - RenderToken
$ws.Range("E44").Value = "  -5.07%  "

# Row 45 - Filecoin
$ws.Range("D45").Value = "4.43"
$ws.Range("E45").Value = "  -3.41%  "

# Row 46 - VeChain
$ws.Range("D46").Value = "0.0309"
$ws.Range("E46").Value = "  -2.93%  "

# Row 47 - dogwifhat
$ws.Range("D47").Value = "2.49"
$ws.Range("E47").Value = "  +8.84%  "

# Row 48 - Bittensor
$ws.Range("D48").Value = "334.56"
$ws.Range("E48").Value = "  +5.25%  "

# Row 49 - ONDO
$ws.Range("E49").Value = "  -2.81%  "

# Row 50 - Cosmos
$ws.Range("E50").Value = "  -2.99%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  -3.05%  "

# Restore the default cell style on column D so no leftover "Text"
# number-format style remains applied (keeps styles.xml aligned with
# the source workbook, which uses the default style for these cells).
$priceRange.Style = "Normal"
